# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Column D (Price) and E (Volume 1h) are stored as text in this sheet, and a
# couple of rows also swap which coin they describe. For D-column values that
# look like plain decimals (e.g. "477.20"), force the cell to Text format
# first so Excel doesn't silently coerce the literal into a Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '55.749.10'
$ws.Range("E2").Value = '  -1.84%  '

# Row 3
$ws.Range("D3").Value = '2.373.65'
$ws.Range("E3").Value = '  -5.05%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '477.20'
$ws.Range("E5").Value = '  -2.45%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.76'
$ws.Range("E6").Value = '  -0.11%  '

# Row 7
$ws.Range("E7").Value = '  +0.21%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.502'
$ws.Range("E8").Value = '  -2.34%  '

# Row 9
$ws.Range("D9").Value = '2.369.09'
$ws.Range("E9").Value = '  -5.99%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0969'
$ws.Range("E10").Value = '  -0.89%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.40'
$ws.Range("E11").Value = '  -6.80%  '

# Row 12
$ws.Range("E12").Value = '  -3.16%  '

# Row 13
$ws.Range("E13").Value = '  +0.72%  '

# Row 14
$ws.Range("D14").Value = '2.788.24'
$ws.Range("E14").Value = '  -5.12%  '

# Row 15
$ws.Range("D15").Value = '55.853.37'
$ws.Range("E15").Value = '  -1.40%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.30'
$ws.Range("E16").Value = '  -4.58%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000132'
$ws.Range("E17").Value = '  -3.77%  '

# Row 18
$ws.Range("D18").Value = '2.382.39'
$ws.Range("E18").Value = '  -5.41%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.57'
$ws.Range("E19").Value = '  +0.97%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '313.88'
$ws.Range("E20").Value = '  -2.34%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.67'
$ws.Range("E21").Value = '  -5.37%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.04%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.65'
$ws.Range("E23").Value = '  -3.23%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '56.61'
$ws.Range("E24").Value = '  -3.85%  '

# Row 25
$ws.Range("E25").Value = '  +0.35%  '

# Row 26
$ws.Range("E26").Value = '  -4.44%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.156'
$ws.Range("E27").Value = '  -6.85%  '

# Row 28
$ws.Range("D28").Value = '2.488.17'
$ws.Range("E28").Value = '  -5.05%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.22'
$ws.Range("E29").Value = '  -5.56%  '

# Row 30
$ws.Range("E30").Value = '  -4.03%  '

# Row 31
$ws.Range("E31").Value = '  +0.12%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '146.45'
$ws.Range("E32").Value = '  -1.68%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.98'
$ws.Range("E33").Value = '  -1.62%  '

# Row 34
$ws.Range("E34").Value = '  -2.06%  '

# Row 35
$ws.Range("E35").Value = '  -3.44%  '

# Row 36
$ws.Range("E36").Value = '  -4.64%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.57'
$ws.Range("E37").Value = '  -4.70%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.833'
$ws.Range("E38").Value = '  -4.25%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '33.40'
$ws.Range("E39").Value = '  -2.51%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  +0.21%  '

# Row 41
$ws.Range("E41").Value = '  -0.26%  '

# Row 42
$ws.Range("E42").Value = '  -5.05%  '

# Row 43
$ws.Range("E43").Value = '  -4.40%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0945'
$ws.Range("E44").Value = '  +3.67%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.576'
$ws.Range("E45").Value = '  -6.91%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.18'
$ws.Range("E46").Value = '  -0.18%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '253.15'
$ws.Range("E47").Value = '  -3.07%  '

# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.57'
$ws.Range("E48").Value = '  -5.33%  '

# Row 49
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0222'
$ws.Range("E49").Value = '  -2.98%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.93'
$ws.Range("E50").Value = '  -4.35%  '

# Row 51
$ws.Range("D51").Value = '1.781.01'
$ws.Range("E51").Value = '  -7.11%  '
